$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 714350.2
$ws.Range("J9").Value = 105.666664
$ws.Range("L9").Value = 105.666664
$ws.Range("N9").Value = -443.666664

$ws.Range("H58").Value = 20834040
$ws.Range("I58").Value = 22727588
$ws.Range("K58").Value = 68182764
$ws.Range("M58").Value = -68182614

$ws.Range("H98").Value = 15430.818
$ws.Range("J98").Value = 980.6667
$ws.Range("L98").Value = 980.6667
$ws.Range("N98").Value = -3976.6667

$ws.Range("H122").Value = 15430.818
$ws.Range("J122").Value = 980.6667
$ws.Range("L122").Value = 2942.0001
$ws.Range("N122").Value = -7842.0001

$ws.Range("H137").Value = 1459.8572
$ws.Range("J137").Value = 2001
$ws.Range("L137").Value = 6003
$ws.Range("N137").Value = -11103

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 72500
$ws.Range("I33").Value = 120000
$ws.Range("J33").Value = 25000
$ws.Range("K33").Value = 120000
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = -119671
$ws.Range("N33").Value = -25658

$ws.Range("H110").Value = 1574
$ws.Range("I110").Value = 570
$ws.Range("K110").Value = 570
$ws.Range("M110").Value = 1475

$ws.Range("H122").Value = 2001.6666
$ws.Range("I122").Value = 1572.25
$ws.Range("K122").Value = 4716.75
$ws.Range("M122").Value = -2266.75

$ws.Range("H132").Value = 3210.818
$ws.Range("I132").Value = 3331.9
$ws.Range("K132").Value = 9995.700000000001
$ws.Range("M132").Value = -7465.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 110000
$ws.Range("J57").Value = 110000
$ws.Range("L57").Value = 110000
$ws.Range("N57").Value = -111440

$ws.Range("H99").Value = 3347.647
$ws.Range("I99").Value = 2001.6666
$ws.Range("J99").Value = 4081.818
$ws.Range("K99").Value = 2001.6666
$ws.Range("L99").Value = 4081.818
$ws.Range("M99").Value = -503.6666
$ws.Range("N99").Value = -7077.818

$ws.Range("H107").Value = 3603.8333
$ws.Range("I107").Value = 1955.1875
$ws.Range("J107").Value = 6901.125
$ws.Range("K107").Value = 1955.1875
$ws.Range("L107").Value = 6901.125
$ws.Range("M107").Value = -35.1875
$ws.Range("N107").Value = -10741.125

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").Value = 0

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").Value = 0

$ws.Range("H134").Value = 1821
$ws.Range("I134").Value = 1665.2759
$ws.Range("K134").Value = 4995.8277
$ws.Range("M134").Value = -2460.8277

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").Value = 0

$ws.Range("H136").Value = 110000
$ws.Range("J136").Value = 110000
$ws.Range("L136").Value = 110000
$ws.Range("N136").Value = -120200

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("N137").Value = 0

$ws.Range("H138").Value = 197142.86
$ws.Range("I138").Value = 60000.332
$ws.Range("K138").Value = 60000.332
$ws.Range("M138").Value = -54860.332

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 969
$ws.Range("I16").Value = 957.6
$ws.Range("J16").Value = 997.5
$ws.Range("K16").Value = 957.6
$ws.Range("L16").Value = 997.5
$ws.Range("M16").Value = -670.6
$ws.Range("N16").Value = -1571.5

$ws.Range("H107").Value = 313.41666
$ws.Range("I107").Value = 319.27274
$ws.Range("J107").Value = 249
$ws.Range("K107").Value = 319.27274
$ws.Range("L107").Value = 249
$ws.Range("M107").Value = 1600.72726
$ws.Range("N107").Value = -4089

$ws.Range("H113").Value = 969
$ws.Range("I113").Value = 957.6
$ws.Range("J113").Value = 997.5
$ws.Range("K113").Value = 957.6
$ws.Range("L113").Value = 997.5
$ws.Range("M113").Value = 1212.4
$ws.Range("N113").Value = -5337.5

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("N127").Value = 0

$ws.Range("H134").Value = 1614.6765
$ws.Range("I134").Value = 1634.9667
$ws.Range("K134").Value = 4904.9001
$ws.Range("M134").Value = -2369.9001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6079268.5
$ws.Range("I4").Value = 1413670.8
$ws.Range("K4").Value = 4241012.4
$ws.Range("M4").Value = -4240900.4

$ws.Range("H68").Value = 993.3125
$ws.Range("J68").Value = 1020.9286
$ws.Range("L68").Value = 3062.7858
$ws.Range("N68").Value = -4684.7858

$ws.Range("H71").Value = 993.3125
$ws.Range("J71").Value = 1020.9286
$ws.Range("L71").Value = 9188.357399999999
$ws.Range("N71").Value = -17300.3574

$ws.Range("H97").Value = 1628.8889
$ws.Range("J97").Value = 1738.75
$ws.Range("L97").Value = 5216.25
$ws.Range("N97").Value = -6208.25

$ws.Range("H107").Value = 1655.3684
$ws.Range("I107").Value = 3078.7144
$ws.Range("J107").Value = 825.0833
$ws.Range("K107").Value = 9236.143199999999
$ws.Range("L107").Value = 2475.2499
$ws.Range("M107").Value = -7316.143199999999
$ws.Range("N107").Value = -6315.2499

$ws.Range("H122").Value = 3948.8333
$ws.Range("J122").Value = 5498.5
$ws.Range("L122").Value = 49486.5
$ws.Range("N122").Value = -54386.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19600
$ws.Range("I70").Value = 19600
$ws.Range("K70").Value = 19600
$ws.Range("M70").Value = -19330

$ws.Range("H73").Value = 19600
$ws.Range("I73").Value = 19600
$ws.Range("K73").Value = 19600
$ws.Range("M73").Value = -18664

$ws.Range("H107").Value = 2168.2307
$ws.Range("I107").Value = 439
$ws.Range("J107").Value = 3249
$ws.Range("K107").Value = 439
$ws.Range("L107").Value = 3249
$ws.Range("M107").Value = 1481
$ws.Range("N107").Value = -7089

$ws.Range("H113").Value = 5823.8066
$ws.Range("I113").Value = 3218.55
$ws.Range("K113").Value = 3218.55
$ws.Range("M113").Value = -1048.55

$ws.Range("H132").Value = 3516.5
$ws.Range("I132").Value = 2034.5
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 6103.5
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -3573.5
$ws.Range("N132").Value = -20055.5

$ws.Range("H140").Value = 163244.5
$ws.Range("J140").Value = 230780
$ws.Range("L140").Value = 230780
$ws.Range("N140").Value = -241140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4532.6562
$ws.Range("I68").Value = 2718
$ws.Range("J68").Value = 5944.0557
$ws.Range("K68").Value = 2718
$ws.Range("L68").Value = 5944.0557
$ws.Range("M68").Value = -1969
$ws.Range("N68").Value = -7442.0557

$ws.Range("H71").Value = 4532.6562
$ws.Range("I71").Value = 2718
$ws.Range("J71").Value = 5944.0557
$ws.Range("K71").Value = 13590
$ws.Range("L71").Value = 29720.2785
$ws.Range("M71").Value = -9846
$ws.Range("N71").Value = -37208.2785

$ws.Range("H93").Value = 4838
$ws.Range("I93").Value = 899.5714
$ws.Range("K93").Value = 899.5714
$ws.Range("M93").Value = 348.4286

$ws.Range("H132").Value = 5531.067
$ws.Range("I132").Value = 7160.625
$ws.Range("J132").Value = 3668.7144
$ws.Range("K132").Value = 21481.875
$ws.Range("L132").Value = 11006.1432
$ws.Range("M132").Value = -18951.875
$ws.Range("N132").Value = -16066.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 11399.4
$ws.Range("J74").Value = 11399.4
$ws.Range("L74").Value = 11399.4
$ws.Range("N74").Value = -13271.4

$ws.Range("H77").Value = 11399.4
$ws.Range("J77").Value = 11399.4
$ws.Range("L77").Value = 34198.2
$ws.Range("N77").Value = -43558.2

$ws.Range("H80").Value = 75000
$ws.Range("J80").Value = 75000
$ws.Range("L80").Value = 75000
$ws.Range("N80").Value = -76996

$ws.Range("H82").Value = 333366660
$ws.Range("J82").Value = 333366660
$ws.Range("L82").Value = 333366660
$ws.Range("N82").Value = -333367426

$ws.Range("H83").Value = 75000
$ws.Range("J83").Value = 75000
$ws.Range("L83").Value = 225000
$ws.Range("N83").Value = -234984

$ws.Range("H85").Value = 333366660
$ws.Range("J85").Value = 333366660
$ws.Range("L85").Value = 333366660
$ws.Range("N85").Value = -333369312

$ws.Range("H96").Value = 3972.6667
$ws.Range("J96").Value = 5500.5
$ws.Range("L96").Value = 5500.5
$ws.Range("N96").Value = -8246.5

$ws.Range("H107").Value = 1540.6207
$ws.Range("I107").Value = 1527.16
$ws.Range("K107").Value = 4581.48
$ws.Range("M107").Value = -2661.48

$ws.Range("H122").Value = 4009.2222
$ws.Range("I122").Value = 2619.25
$ws.Range("J122").Value = 7980.5713
$ws.Range("K122").Value = 7857.75
$ws.Range("L122").Value = 23941.7139
$ws.Range("M122").Value = -5407.75
$ws.Range("N122").Value = -28841.7139
